# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (col I) and DialogAct (col J) values for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 14; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 15; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 34; I = '%'; J = 'Uninterpretable' }
    @{ Row = 39; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 62; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 64; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 78; I = 'ba'; J = 'Appreciation' }
    @{ Row = 79; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 85; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 102; I = '%'; J = 'Uninterpretable' }
    @{ Row = 103; I = '%'; J = 'Uninterpretable' }
    @{ Row = 107; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 113; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 126; I = 'ba'; J = 'Appreciation' }
    @{ Row = 127; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 131; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 137; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 138; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 142; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 145; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 149; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 158; I = 'ba'; J = 'Appreciation' }
    @{ Row = 159; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 160; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 161; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 163; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 166; I = 'ba'; J = 'Appreciation' }
    @{ Row = 173; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 174; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 184; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 192; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 194; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 198; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 205; I = 'ba'; J = 'Appreciation' }
    @{ Row = 211; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 212; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 213; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 216; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 219; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 220; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 226; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 230; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 232; I = '%'; J = 'Uninterpretable' }
    @{ Row = 234; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 235; I = 'ba'; J = 'Appreciation' }
    @{ Row = 236; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 237; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 245; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 266; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 278; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 280; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 288; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 293; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 296; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 318; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 324; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 341; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 342; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 347; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 358; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 361; I = 'ba'; J = 'Appreciation' }
    @{ Row = 388; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 395; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 396; I = 'aa'; J = 'Agree/Accept' }
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg.Row, 9).Value = $chg.I
    $ws.Cells.Item($chg.Row, 10).Value = $chg.J
}
